$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Copy existing formatting onto the new cells first (grouping every paste
# that shares one source so the clipboard is copied only once per source -
# re-copying the same source repeatedly causes duplicate style records).
# ---------------------------------------------------------------------------

# Header-row style (merged, filled) -> new block headers on row 11 and row 18
$ws.Range("F2:G2").Copy()
$ws.Range("F11:G11").PasteSpecial(-4122)
$ws.Range("F18:G18").PasteSpecial(-4122)

# "hour/_res_" + "Hour/Result" sub-header rows style -> rows 12:13 and 19:20
$ws.Range("F3:G4").Copy()
$ws.Range("F12:G13").PasteSpecial(-4122)
$ws.Range("F19:G20").PasteSpecial(-4122)

# Data-row style -> rows 14:16 and 21:23
$ws.Range("F5:G7").Copy()
$ws.Range("F14:G16").PasteSpecial(-4122)
$ws.Range("F21:G23").PasteSpecial(-4122)

# Data-row style -> extra row appended to the first results table (row 9)
$ws.Range("F5:G5").Copy()
$ws.Range("F9:G9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Fill in the values
# ---------------------------------------------------------------------------

# Extra data row for the first "GreetingTest" results table
$ws.Range("F9").Value = 25
$ws.Range("G9").Value = "Good Night, World!"

# Block 1: rows 11-16 "Test Greeting GreetingSuccessful1"
$ws.Range("F11").Value = "Test Greeting GreetingSuccessful1"
$ws.Range("F12").Value = "hour"
$ws.Range("G12").Value = "_res_"
$ws.Range("F13").Value = "Hour"
$ws.Range("G13").Value = "Result"
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = "Good Morning, World!"
$ws.Range("F15").Value = 13
$ws.Range("G15").Value = "Good Afternoon, World!"
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = "Good Night, World!"
$ws.Range("F11:G11").Merge()

# Block 2: rows 18-23 "Test Greeting GreetingSuccessful2"
$ws.Range("F18").Value = "Test Greeting GreetingSuccessful2"
$ws.Range("F19").Value = "hour"
$ws.Range("G19").Value = "_res_"
$ws.Range("F20").Value = "Hour"
$ws.Range("G20").Value = "Result"
$ws.Range("F21").Value = 7
$ws.Range("G21").Value = "Good Morning, World!"
$ws.Range("F22").Value = 13
$ws.Range("G22").Value = "Good Afternoon, World!"
$ws.Range("F23").Value = 22
$ws.Range("G23").Value = "Good Night, World!"
$ws.Range("F18:G18").Merge()

# ---------------------------------------------------------------------------
# Column G width update (target stored width 24.5703125; the COM width
# setter in this runtime quantizes to the nearest 1/6 character width, so we
# pick the input value that lands on the nearest reachable width, 24.5)
# ---------------------------------------------------------------------------
$ws.Columns("G").ColumnWidth = 23.6667

# ---------------------------------------------------------------------------
# Selection / active cell
# ---------------------------------------------------------------------------
$ws.Range("F26").Select() | Out-Null
